$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C (LOAN_NUMBER) to fit the longer numeric-text loan numbers
$ws.Columns.Item(3).ColumnWidth = 28.29

# The loan numbers are large integers that must be stored as text (so they
# aren't mangled by floating point/scientific notation) - format column C as Text
$ws.Range("C1:C4").NumberFormat = "@"

# Replace the old loan-number values with the new WO numbers (entered as text)
$ws.Range("C2").Value = "1950718000000550"
$ws.Range("C3").Value = "1950718000000551"
$ws.Range("C4").Value = "1950718000000552"

# Update the outstanding/write-off flags in column E
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 0

# Switch the page to portrait orientation for printing
$ws.PageSetup.Orientation = 1

# Leave the cursor where the user last clicked while reviewing the sheet
$ws.Range("F12").Select()
